# Applies the Nov 26 2023 cryptos price/volume refresh to Sheet1 (rows 2-51).
# Price (col D) and Volume(1h) (col E) are plain text cells (e.g. "37.242.89",
# "  -1.61%  "). A leading apostrophe forces text entry for values that would
# otherwise be auto-parsed as numbers by Excel's COM layer (matches how a
# human typing these into Excel would keep them as text), while leaving the
# cell NumberFormat at General, same as the source file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.242.89"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").Value = "2.046.78"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'230.46"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'56.98"
$ws.Range("E8").Value = "  -4.42%  "

$ws.Range("E9").Value = "  -2.82%  "

$ws.Range("D10").Value = "'0.0767"
$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "'14.65"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").Value = "2.351.73"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "'20.58"
$ws.Range("E14").Value = "  -3.25%  "

$ws.Range("E15").Value = "  -2.75%  "

$ws.Range("E16").Value = "  -2.23%  "

$ws.Range("D17").Value = "2.053.19"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "37.242.61"
$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").Value = "'5.99"
$ws.Range("E19").Value = "  -2.97%  "

$ws.Range("D20").Value = "'69.64"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -3.67%  "

$ws.Range("D22").Value = "'226.67"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("E25").Value = "  -4.07%  "

$ws.Range("D26").Value = "'9.60"
$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").Value = "'169.71"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("E28").Value = "  -4.09%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("E30").Value = "  -5.94%  "

$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("E32").Value = "  -4.54%  "

$ws.Range("E33").Value = "  -2.06%  "

$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("D35").Value = "'2.49"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("E37").Value = "  -4.08%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("D40").Value = "'0.0224"
$ws.Range("E40").Value = "  +3.56%  "

$ws.Range("D41").Value = "'98.36"
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D42").Value = "1.486.24"
$ws.Range("E42").Value = "  +2.77%  "

$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("E44").Value = "  -3.79%  "

$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").Value = "'16.43"
$ws.Range("E46").Value = "  -1.94%  "

$ws.Range("E47").Value = "  -3.69%  "

$ws.Range("D48").Value = "'3.95"
$ws.Range("E48").Value = "  -6.19%  "

$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").Value = "'2.94"
$ws.Range("E50").Value = "  -2.44%  "

$ws.Range("D51").Value = "2.237.43"
$ws.Range("E51").Value = "  -1.42%  "
